$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: new entry - "Remove Duplicates from Sorted Array" (LeetCode problem)
# Sl no.
$ws.Range("A12").Value = 18

# File/title column with a hyperlink to the LeetCode problem page.
$ws.Range("B12").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B12"), "https://leetcode.com/problems/remove-duplicates-from-sorted-array/", "", "", "https://leetcode.com/problems/remove-duplicates-from-sorted-array/") | Out-Null
$ws.Range("B12").Value = "Remove Duplicates from Sorted Array"

# Date column - copy the date format used by the surrounding rows (D10/D11)
# and set the serial date value for 2025-07-20.
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("D12").Value = 45858

# Move the active selection to the newly-edited cell.
$ws.Range("B12").Select() | Out-Null
